$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) — update "想去人数" (wish-to-go count) values in column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 15440
$wsExpo.Range("F9").Value = 15365
$wsExpo.Range("F11").Value = 8959
$wsExpo.Range("F15").Value = 82
$wsExpo.Range("F32").Value = 53
$wsExpo.Range("F34").Value = 243
$wsExpo.Range("F36").Value = 442
$wsExpo.Range("F38").Value = 5483

# Sheet "全部类型" (All Types) — same updates, rows shifted by 2 for the later entries
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 15440
$wsAll.Range("F9").Value = 15365
$wsAll.Range("F11").Value = 8959
$wsAll.Range("F15").Value = 82
$wsAll.Range("F34").Value = 53
$wsAll.Range("F36").Value = 243
$wsAll.Range("F38").Value = 442
$wsAll.Range("F40").Value = 5483
